$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: reorder C/D/E (Qty, Cost, Description) -> (Description, Qty, Cost) ---
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "Qty"
$ws.Range("E1").Value = "Cost"

$ws.Range("C2").Value = "Raspberry Pi Pico"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 4
$ws.Range("F2").Formula = "=D2*E2"
$ws.Range("C3").Value = "Power FETs"
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = 1.9
$ws.Range("F3").Formula = "=D3*E3"
$ws.Range("C4").Value = "Gate driver"
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 1.26
$ws.Range("F4").Formula = "=D4*E4"
$ws.Range("C5").Value = "Diode for gate driver"
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 0.25
$ws.Range("F5").Formula = "=D5*E5"
$ws.Range("C6").Value = "Current shunt amplifier"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 0.57999999999999996
$ws.Range("F6").Formula = "=D6*E6"
$ws.Range("C7").Value = "Current shunt resistor"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0.7
$ws.Range("F7").Formula = "=D7*E7"
$ws.Range("C8").Value = "Charge pump IC"
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 1.1100000000000001
$ws.Range("F8").Formula = "=D8*E8"
$ws.Range("C9").Value = "5V buck regulator"
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 2.89
$ws.Range("F9").Formula = "=D9*E9"
$ws.Range("C10").Value = "220uF Capacitor"
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 1.1399999999999999
$ws.Range("F10").Formula = "=D10*E10"
$ws.Range("C11").Value = "TVS Diode"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 0.47
$ws.Range("F11").Formula = "=D11*E11"
$ws.Range("C12").Value = "47k Resistor"
$ws.Range("D12").Value = 5
$ws.Range("E12").Value = 0.03
$ws.Range("F12").Formula = "=D12*E12"
$ws.Range("C13").Value = "2.2k Resistor"
$ws.Range("D13").Value = 7
$ws.Range("E13").Value = 0.03
$ws.Range("F13").Formula = "=D13*E13"
$ws.Range("C14").Value = "22nF Capacitor"
$ws.Range("D14").Value = 4
$ws.Range("E14").Value = 0.08
$ws.Range("F14").Formula = "=D14*E14"
$ws.Range("C15").Value = "1uF Capacitor"
$ws.Range("D15").Value = 13
$ws.Range("E15").Value = 0.19
$ws.Range("F15").Formula = "=D15*E15"
$ws.Range("C16").Value = "Fuse holder"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 1.33
$ws.Range("F16").Formula = "=D16*E16"
$ws.Range("C17").Value = "Reset pushbutton"
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0.13
$ws.Range("F17").Formula = "=D17*E17"
$ws.Range("C18").Value = "4 pin female VAL-U-LOK"
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0.68
$ws.Range("F18").Formula = "=D18*E18"
$ws.Range("C19").Value = "6 pin female VAL-U-LOK"
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 0.65
$ws.Range("F19").Formula = "=D19*E19"
$ws.Range("C20").Value = "4 pin male VAL-U-LOK"
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0.21
$ws.Range("F20").Formula = "=D20*E20"
$ws.Range("C21").Value = "6 pin male VAL-U-LOK"
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 0.24
$ws.Range("F21").Formula = "=D21*E21"
$ws.Range("C22").Value = "Crimp pin VAL-U-LOK"
$ws.Range("D22").Value = 10
$ws.Range("E22").Value = 0.05
$ws.Range("F22").Formula = "=D22*E22"
$ws.Range("C23").Value = "30A fuse"
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0.33
$ws.Range("F23").Formula = "=D23*E23"

# --- sheetPr / page setup: enable "fit to page" and set scale + landscape orientation ---
$ws.PageSetup.Zoom = 87
$ws.PageSetup.Orientation = 2
$ws.PageSetup.FitToPagesTall = 1

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 25.42
$ws.Columns.Item(2).ColumnWidth = 19.92
$ws.Columns.Item(3).ColumnWidth = 19.92
$ws.Columns.Item(6).ColumnWidth = 7.42

# --- Selection ---
$ws.Range("G11").Select()
